$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "Ticket_Category" column (F) ----------------------------
$ws.Range("F1").Value = "Ticket_Category"

# Give it the same look as the other header cells (bold 10pt Segoe UI, blue)
# by copying E1's formatting only - avoids creating extra intermediate styles.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Column F width (matches the bestFit width used by the other header columns)
$ws.Columns("F").ColumnWidth = 16.9

# --- Update the AutoFilter so it now spans A1:F1 -------------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:F1").AutoFilter()

# Keep the workbook-level _FilterDatabase defined name in sync with the filter range
$fdb = $wb.Names.Item("IncidentReport!_FilterDatabase")
$fdb.RefersTo = "=IncidentReport!`$A`$1:`$F`$1"

# --- Update the active selection -----------------------------------------
$ws.Range("A2").Select()
